$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1956.6666
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 1956.6666
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 5869.9998
$ws.Range("N17").Value = -6205.9998
$ws.Range("M17").ClearContents()

$ws.Range("H29").Value = 2173.6667
$ws.Range("J29").Value = 3234.75
$ws.Range("L29").Value = 9704.25
$ws.Range("N29").Value = -10266.25

$ws.Range("H43").Value = 6999.857
$ws.Range("I43").Value = 4800
$ws.Range("K43").Value = 4800
$ws.Range("M43").Value = -4731

$ws.Range("H86").Value = 367001.34
$ws.Range("I86").Value = 1000
$ws.Range("K86").Value = 1000
$ws.Range("M86").Value = 123

$ws.Range("H88").Value = 1914
$ws.Range("I88").Value = 1769.8
$ws.Range("J88").Value = 1969.4615
$ws.Range("K88").Value = 1769.8
$ws.Range("L88").Value = 1969.4615
$ws.Range("M88").Value = -1363.8
$ws.Range("N88").Value = -2781.4615

$ws.Range("H89").Value = 367001.34
$ws.Range("I89").Value = 1000
$ws.Range("K89").Value = 5000
$ws.Range("M89").Value = 616

$ws.Range("H91").Value = 1914
$ws.Range("I91").Value = 1769.8
$ws.Range("J91").Value = 1969.4615
$ws.Range("K91").Value = 1769.8
$ws.Range("L91").Value = 1969.4615
$ws.Range("M91").Value = -365.8
$ws.Range("N91").Value = -4777.461499999999

$ws.Range("H100").Value = 875.5714
$ws.Range("I100").Value = 896.5
$ws.Range("K100").Value = 896.5
$ws.Range("M100").Value = -355.5

$ws.Range("H113").Value = 8756.857
$ws.Range("J113").Value = 8999.5
$ws.Range("L113").Value = 8999.5
$ws.Range("N113").Value = -15507.5

$ws.Range("H137").Value = 1232.9
$ws.Range("I137").Value = 865.6087
$ws.Range("J137").Value = 2439.7144
$ws.Range("K137").Value = 2596.8261
$ws.Range("L137").Value = 7319.1432
$ws.Range("M137").Value = -46.82610000000022
$ws.Range("N137").Value = -12419.1432

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1901.7142
$ws.Range("I45").Value = 1377.25
$ws.Range("K45").Value = 1377.25
$ws.Range("M45").Value = -1000.25

$ws.Range("H122").Value = 12323.296
$ws.Range("I122").Value = 9113.708000000001
$ws.Range("K122").Value = 27341.124
$ws.Range("M122").Value = -24891.124

$ws.Range("H132").Value = 5964.8125
$ws.Range("I132").Value = 5560.5557
$ws.Range("J132").Value = 6484.5713
$ws.Range("K132").Value = 16681.6671
$ws.Range("L132").Value = 19453.7139
$ws.Range("M132").Value = -14151.6671
$ws.Range("N132").Value = -24513.7139

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H95").Value = 16463.572
$ws.Range("I95").Value = 3999
$ws.Range("J95").Value = 18541
$ws.Range("K95").Value = 3999
$ws.Range("L95").Value = 18541
$ws.Range("M95").Value = -1253
$ws.Range("N95").Value = -24033

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 3291
$ws.Range("I16").Value = 3618.2856
$ws.Range("K16").Value = 3618.2856
$ws.Range("M16").Value = -3331.2856

$ws.Range("H22").Value = 74877.55499999999
$ws.Range("J22").Value = 20856.572
$ws.Range("L22").Value = 20856.572
$ws.Range("N22").Value = -21556.572

$ws.Range("H31").Value = 2350.923
$ws.Range("I31").Value = 1597.2273
$ws.Range("K31").Value = 1597.2273
$ws.Range("M31").Value = -1302.2273

$ws.Range("H34").Value = 2350.923
$ws.Range("I34").Value = 1597.2273
$ws.Range("K34").Value = 1597.2273
$ws.Range("M34").Value = -1395.2273

$ws.Range("H47").Value = 17000
$ws.Range("J47").Value = 17000
$ws.Range("L47").Value = 17000
$ws.Range("N47").Value = -18132

$ws.Range("H113").Value = 3291
$ws.Range("I113").Value = 3618.2856
$ws.Range("K113").Value = 3618.2856
$ws.Range("M113").Value = -1448.2856

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 98000
$ws.Range("J37").Value = 98000
$ws.Range("L37").Value = 294000
$ws.Range("N37").Value = -294224

$ws.Range("H97").Value = 1756.0769
$ws.Range("I97").Value = 979.6
$ws.Range("J97").Value = 2241.375
$ws.Range("K97").Value = 2938.8
$ws.Range("L97").Value = 6724.125
$ws.Range("M97").Value = -2442.8
$ws.Range("N97").Value = -7716.125

$ws.Range("H131").Value = 1419.6086

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H47").Value = 26638.666
$ws.Range("J47").Value = 30001
$ws.Range("L47").Value = 30001
$ws.Range("N47").Value = -31137

$ws.Range("H80").Value = 4913.5
$ws.Range("J80").Value = 4750
$ws.Range("L80").Value = 4750
$ws.Range("N80").Value = -6746

$ws.Range("H83").Value = 4913.5
$ws.Range("J83").Value = 4750
$ws.Range("L83").Value = 23750
$ws.Range("N83").Value = -33734

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3431.9614
$ws.Range("I40").Value = 3431.9614
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 3431.9614
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = -3295.9614
$ws.Range("N40").ClearContents()

$ws.Range("H46").Value = 4287.5293
$ws.Range("J46").Value = 6126.857
$ws.Range("L46").Value = 6126.857
$ws.Range("N46").Value = -6502.857

$ws.Range("H82").Value = 127486.625
$ws.Range("I82").Value = 1974.5
$ws.Range("J82").Value = 169324
$ws.Range("K82").Value = 1974.5
$ws.Range("L82").Value = 169324
$ws.Range("M82").Value = -1613.5
$ws.Range("N82").Value = -170046

$ws.Range("H85").Value = 127486.625
$ws.Range("I85").Value = 1974.5
$ws.Range("J85").Value = 169324
$ws.Range("K85").Value = 1974.5
$ws.Range("L85").Value = 169324
$ws.Range("M85").Value = -726.5
$ws.Range("N85").Value = -171820

$ws.Range("H122").Value = 6270.154
$ws.Range("I122").Value = 5107.2
$ws.Range("J122").Value = 6997
$ws.Range("K122").Value = 15321.6
$ws.Range("L122").Value = 20991
$ws.Range("M122").Value = -12871.6
$ws.Range("N122").Value = -25891

$ws.Range("H132").Value = 57140.61
$ws.Range("I132").Value = 60295.94
$ws.Range("K132").Value = 180887.82
$ws.Range("M132").Value = -178357.82

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H9").Value = 0
$ws.Range("I9").Value = 0
$ws.Range("J9").Value = 0
$ws.Range("K9").Value = 0
$ws.Range("L9").Value = 0
$ws.Range("M9").ClearContents()
$ws.Range("N9").ClearContents()

$ws.Range("H14").Value = 18597.8
$ws.Range("J14").Value = 17997.25
$ws.Range("L14").Value = 17997.25
$ws.Range("N14").Value = -18333.25

$ws.Range("H21").Value = 26400
$ws.Range("J21").Value = 26400
$ws.Range("L21").Value = 26400
$ws.Range("N21").Value = -26870

$ws.Range("H33").Value = 36840.332
$ws.Range("I33").Value = 29000
$ws.Range("J33").Value = 40760.5
$ws.Range("K33").Value = 29000
$ws.Range("L33").Value = 40760.5
$ws.Range("M33").Value = -28750
$ws.Range("N33").Value = -41260.5

$ws.Range("H35").Value = 26400
$ws.Range("J35").Value = 26400
$ws.Range("L35").Value = 26400
$ws.Range("N35").Value = -26980

$ws.Range("H36").Value = 36840.332
$ws.Range("I36").Value = 29000
$ws.Range("J36").Value = 40760.5
$ws.Range("K36").Value = 29000
$ws.Range("L36").Value = 40760.5
$ws.Range("M36").Value = -28750
$ws.Range("N36").Value = -41260.5

$ws.Range("H40").Value = 35000
$ws.Range("I40").Value = 35000
$ws.Range("J40").Value = 35000
$ws.Range("K40").Value = 35000
$ws.Range("L40").Value = 35000
$ws.Range("M40").Value = -34851
$ws.Range("N40").Value = -35298

$ws.Range("H47").Value = 16799
$ws.Range("J47").Value = 16799
$ws.Range("L47").Value = 16799
$ws.Range("N47").Value = -17943

$ws.Range("H100").Value = 1373.0667
$ws.Range("I100").Value = 1543
$ws.Range("J100").Value = 1178.8572
$ws.Range("K100").Value = 3086
$ws.Range("L100").Value = 2357.7144
$ws.Range("M100").Value = -2545
$ws.Range("N100").Value = -3439.7144

$ws.Range("H117").Value = 95204.5
$ws.Range("J117").Value = 95204.5
$ws.Range("L117").Value = 95204.5
$ws.Range("N117").Value = -104382.5

$ws.Range("H122").Value = 3045.5789
$ws.Range("I122").Value = 2705.0715
$ws.Range("J122").Value = 3999
$ws.Range("K122").Value = 8115.2145
$ws.Range("L122").Value = 11997
$ws.Range("M122").Value = -5665.2145
$ws.Range("N122").Value = -16897

$ws.Range("I126").Value = 1244.6666
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 3733.9998
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -1263.9998
$ws.Range("N126").ClearContents()

$ws.Range("H136").Value = 5215
$ws.Range("I136").Value = 5039.6
$ws.Range("K136").Value = 15118.8
$ws.Range("M136").Value = -12568.8
